# Mise à jour du TPI 16.05.2023
# Update "Planning effectif" sheet: fill in actual time-spent values (col B)
# for several tasks, plus the 10th day (col L) entries on the last three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planning effectif")

# Column B ("Temps nécessaire") -- set actual elapsed time for tasks that
# had no estimate recorded yet (now equal to the task's realised total).
$ws.Range("B6").Value = 0.14583333333333334
$ws.Range("B8").Value = 0.20138888888888887
$ws.Range("B9").Value = 0.06597222222222222
$ws.Range("B10").Value = 0.0763888888888889
$ws.Range("B11").Value = 0.006944444444444444
$ws.Range("B12").Value = 0.16319444444444445
$ws.Range("B15").Value = 0.10069444444444443
$ws.Range("B21").Value = 0.22569444444444445
$ws.Range("B22").Value = 0.19444444444444445

# Column L ("10e jour") -- new time entries logged on day 10.
$ws.Range("L24").Value = 0.0763888888888889
$ws.Range("L25").Value = 0.013888888888888888
$ws.Range("L26").Value = 0.24305555555555555

# Update the active selection to reflect where the user last clicked.
$ws.Activate()
$ws.Range("L27").Select()
